# Applies the cryptos price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.888.03"
$ws.Range("E2").Value = "  -4.26%  "

# Row 3
$ws.Range("D3").Value = "3.496.58"
$ws.Range("E3").Value = "  -4.15%  "

# Row 4
$style = $ws.Range("D4").Style
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'576.32"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -2.16%  "

# Row 6
$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'169.98"
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = "  -6.18%  "

# Row 7
$style = $ws.Range("D7").Style
$ws.Range("D7").Value = "'0.614"
$ws.Range("D7").Style = $style
$ws.Range("E7").Value = "  -0.49%  "

# Row 8
$ws.Range("D8").Value = "3.473.72"
$ws.Range("E8").Value = "  -4.67%  "

# Row 9
$ws.Range("E9").Value = "  +0.39%  "

# Row 10
$style = $ws.Range("D10").Style
$ws.Range("D10").Value = "'0.188"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -7.66%  "

# Row 11
$style = $ws.Range("D11").Style
$ws.Range("D11").Value = "'6.68"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +12.81%  "

# Row 12
$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.595"
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = "  -2.41%  "

# Row 13
$style = $ws.Range("D13").Style
$ws.Range("D13").Value = "'46.91"
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = "  -6.08%  "

# Row 14
$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'0.0000273"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  -4.84%  "

# Row 15
$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'685.98"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  +0.06%  "

# Row 16
$ws.Range("D16").Value = "4.061.18"
$ws.Range("E16").Value = "  -4.05%  "

# Row 17
$style = $ws.Range("D17").Style
$ws.Range("D17").Value = "'8.68"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  -4.19%  "

# Row 18
$ws.Range("D18").Value = "68.920.74"
$ws.Range("E18").Value = "  -4.17%  "

# Row 19
$ws.Range("D19").Value = "3.502.77"
$ws.Range("E19").Value = "  -4.48%  "

# Row 20
$ws.Range("E20").Value = "  -1.93%  "

# Row 21
$style = $ws.Range("D21").Style
$ws.Range("D21").Value = "'17.32"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  -5.77%  "

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").Value = "'11.08"
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = "  -5.06%  "

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").Value = "'0.905"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -3.93%  "

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'16.44"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -8.00%  "

# Row 25
$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'97.09"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  -6.16%  "

# Row 26
$style = $ws.Range("D26").Style
$ws.Range("D26").Value = "'3.82"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  -5.09%  "

# Row 27
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$style = $ws.Range("D27").Style
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  -0.15%  "

# Row 28
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$style = $ws.Range("D28").Style
$ws.Range("D28").Value = "'2.65"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  -7.40%  "

# Row 29
$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'9.38"
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = "  -6.44%  "

# Row 30
$style = $ws.Range("D30").Style
$ws.Range("D30").Value = "'32.99"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  -6.16%  "

# Row 31
$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'8.79"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  -5.26%  "

# Row 32
$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'3.16"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  -7.93%  "

# Row 33
$style = $ws.Range("D33").Style
$ws.Range("D33").Value = "'7.27"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  -1.21%  "

# Row 34
$ws.Range("E34").Value = "  -7.26%  "

# Row 35
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'570.62"
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = "  -3.70%  "

# Row 36
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'3.71"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  -11.24%  "

# Row 37
$style = $ws.Range("D37").Style
$ws.Range("D37").Value = "'10.80"
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = "  -4.94%  "

# Row 38
$style = $ws.Range("D38").Style
$ws.Range("D38").Value = "'0.104"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  -4.38%  "

# Row 39
$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'57.22"
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = "  -3.90%  "

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  +0.13%  "

# Row 41
$ws.Range("D41").Value = "3.472.76"
$ws.Range("E41").Value = "  -5.91%  "

# Row 42
$style = $ws.Range("D42").Style
$ws.Range("D42").Value = "'0.0438"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  -7.30%  "

# Row 43
$style = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.136"
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = "  -4.73%  "

# Row 44
$style = $ws.Range("D44").Style
$ws.Range("D44").Value = "'0.335"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  -3.74%  "

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'33.08"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  -7.58%  "

# Row 46
$ws.Range("D46").Value = "0.0₃0698"
$ws.Range("E46").Value = "  -9.16%  "

# Row 47
$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'2.88"
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = "  +2.20%  "

# Row 48
$style = $ws.Range("D48").Style
$ws.Range("D48").Value = "'2.57"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  -7.41%  "

# Row 49
$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'0.131"
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = "  -1.93%  "

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'133.57"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +0.96%  "

# Row 51
$style = $ws.Range("D51").Style
$ws.Range("D51").Value = "'0.148"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  -1.72%  "
